# Updated cryptos list on Tue May 23 03:07:46 UTC 2023 with GitHub Actions
#
# Refreshes the per-coin Price (col D) and Volume(1h) (col E) snapshot values
# pulled from coinranking.com, and fixes the HuobiToken/ImmutableX row order
# (rows 32-33 swapped places, each keeping its own Coin/Link/Price/Volume).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price strings such as "314.41" read as numbers through normal COM assignment;
# Excel keeps them as literal text (matching the source workbook) when the entry
# is apostrophe-prefixed, exactly like typing `'314.41` into the cell by hand.
function Set-CellText($range, [string]$text) {
    $range.Value = "'" + $text
}

# Row 2
$ws.Range("D2").Value = '27.235.07'
$ws.Range("E2").Value = '  +1.48%  '

# Row 3
$ws.Range("D3").Value = '1.849.46'
$ws.Range("E3").Value = '  +1.89%  '

# Row 4
$ws.Range("E4").Value = '  -0.36%  '

# Row 5
Set-CellText $ws.Range("D5") '314.41'
$ws.Range("E5").Value = '  +1.99%  '

# Row 6
$ws.Range("E6").Value = '  -0.37%  '

# Row 7
Set-CellText $ws.Range("D7") '0.4646'
$ws.Range("E7").Value = '  +0.54%  '

# Row 8
Set-CellText $ws.Range("D8") '0.3708'
$ws.Range("E8").Value = '  +1.79%  '

# Row 9
Set-CellText $ws.Range("D9") '0.07373'
$ws.Range("E9").Value = '  +2.15%  '

# Row 10
Set-CellText $ws.Range("D10") '0.8859'
$ws.Range("E10").Value = '  +3.46%  '

# Row 11
Set-CellText $ws.Range("D11") '0.07917'
$ws.Range("E11").Value = '  +5.17%  '

# Row 12
Set-CellText $ws.Range("D12") '20.01'
$ws.Range("E12").Value = '  +1.63%  '

# Row 13
$ws.Range("D13").Value = '1.898.77'
$ws.Range("E13").Value = '  +3.55%  '

# Row 14
Set-CellText $ws.Range("D14") '5.386'
$ws.Range("E14").Value = '  +1.14%  '

# Row 15
Set-CellText $ws.Range("D15") '6.597'
$ws.Range("E15").Value = '  +1.34%  '

# Row 16
Set-CellText $ws.Range("D16") '92.13'
$ws.Range("E16").Value = '  +0.52%  '

# Row 17
$ws.Range("E17").Value = '  -0.25%  '

# Row 18
Set-CellText $ws.Range("D18") '0.000008928'
$ws.Range("E18").Value = '  +3.96%  '

# Row 19
Set-CellText $ws.Range("D19") '1.004'
$ws.Range("E19").Value = '  -0.40%  '

# Row 20
Set-CellText $ws.Range("D20") '14.88'
$ws.Range("E20").Value = '  +3.24%  '

# Row 21
$ws.Range("D21").Value = '27.259.58'
$ws.Range("E21").Value = '  +1.03%  '

# Row 22
Set-CellText $ws.Range("D22") '5.142'
$ws.Range("E22").Value = '  -0.04%  '

# Row 23
Set-CellText $ws.Range("D23") '10.58'
$ws.Range("E23").Value = '  +0.75%  '

# Row 24
$ws.Range("D24").Value = '2.074.50'
$ws.Range("E24").Value = '  +0.46%  '

# Row 25
Set-CellText $ws.Range("D25") '152.89'
$ws.Range("E25").Value = '  +1.02%  '

# Row 26
Set-CellText $ws.Range("D26") '1.862'
$ws.Range("E26").Value = '  +0.81%  '

# Row 27
Set-CellText $ws.Range("D27") '18.51'
$ws.Range("E27").Value = '  +2.26%  '

# Row 28
Set-CellText $ws.Range("D28") '2.073'
$ws.Range("E28").Value = '  +0.56%  '

# Row 29
Set-CellText $ws.Range("D29") '5.147'
$ws.Range("E29").Value = '  +1.37%  '

# Row 30
Set-CellText $ws.Range("D30") '117.18'
$ws.Range("E30").Value = '  +1.98%  '

# Row 31
Set-CellText $ws.Range("D31") '0.08893'
$ws.Range("E31").Value = '  +0.42%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellText $ws.Range("D32") '0.7470'
$ws.Range("E32").Value = '  +4.24%  '

# Row 33
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-CellText $ws.Range("D33") '2.974'
$ws.Range("E33").Value = '  +0.77%  '

# Row 34
Set-CellText $ws.Range("D34") '4.477'
$ws.Range("E34").Value = '  +1.84%  '

# Row 35
Set-CellText $ws.Range("D35") '1.145'
$ws.Range("E35").Value = '  +1.34%  '

# Row 36
Set-CellText $ws.Range("D36") '2.554'
$ws.Range("E36").Value = '  +5.38%  '

# Row 37
Set-CellText $ws.Range("D37") '1.081'
$ws.Range("E37").Value = '  +0.83%  '

# Row 38
Set-CellText $ws.Range("D38") '0.05283'
$ws.Range("E38").Value = '  +0.99%  '

# Row 39
Set-CellText $ws.Range("D39") '0.01954'
$ws.Range("E39").Value = '  +2.17%  '

# Row 40
Set-CellText $ws.Range("D40") '2.973'
$ws.Range("E40").Value = '  +1.86%  '

# Row 41
Set-CellText $ws.Range("D41") '7.110'
$ws.Range("E41").Value = '  -0.48%  '

# Row 42
Set-CellText $ws.Range("D42") '0.5185'
$ws.Range("E42").Value = '  +1.14%  '

# Row 43
$ws.Range("E43").Value = '  +1.08%  '

# Row 44
Set-CellText $ws.Range("D44") '8.331'
$ws.Range("E44").Value = '  +2.06%  '

# Row 45
Set-CellText $ws.Range("D45") '0.4873'
$ws.Range("E45").Value = '  +1.77%  '

# Row 46
Set-CellText $ws.Range("D46") '10.23'
$ws.Range("E46").Value = '  +1.01%  '

# Row 47
Set-CellText $ws.Range("D47") '1.004'
$ws.Range("E47").Value = '  -0.33%  '

# Row 48
Set-CellText $ws.Range("D48") '103.07'
$ws.Range("E48").Value = '  +0.17%  '

# Row 49
Set-CellText $ws.Range("D49") '1.635'
$ws.Range("E49").Value = '  +1.41%  '

# Row 50
Set-CellText $ws.Range("D50") '0.06235'
$ws.Range("E50").Value = '  +0.63%  '

# Row 51
Set-CellText $ws.Range("D51") '65.58'
$ws.Range("E51").Value = '  +2.37%  '
